$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.8251781463623
$ws.Range("B3").Value = 18.17774391174316
$ws.Range("B4").Value = 17.91836929321289
$ws.Range("B5").Value = 18.65038871765137
$ws.Range("B6").Value = 21.23887062072754
$ws.Range("B7").Value = 21.25138664245605
$ws.Range("B8").Value = 24.5553035736084
$ws.Range("B9").Value = 26.6196117401123
$ws.Range("B10").Value = 24.83595848083496
$ws.Range("B11").Value = 24.18892860412598
$ws.Range("B12").Value = 15.07848834991455
$ws.Range("B13").Value = 13.82089996337891
$ws.Range("B14").Value = 12.58361339569092
$ws.Range("B15").Value = 13.29953384399414
$ws.Range("B16").Value = 12.66551303863525
$ws.Range("B17").Value = 12.92565250396729
$ws.Range("B18").Value = 22.69183349609375
$ws.Range("B19").Value = 24.6817512512207
$ws.Range("B20").Value = 23.97696685791016
$ws.Range("B21").Value = 32.19055557250977
$ws.Range("B22").Value = 32.30900955200195
$ws.Range("B23").Value = 31.03658485412598
$ws.Range("B24").Value = 20.4122257232666
$ws.Range("B25").Value = 16.83600044250488
$ws.Range("B26").Value = 15.17129421234131
$ws.Range("B27").Value = 12.56902122497559
$ws.Range("B28").Value = 11.32696437835693
$ws.Range("B29").Value = 24.93119239807129
$ws.Range("B30").Value = 36.83232879638672
$ws.Range("B31").Value = 34.52827072143555
$ws.Range("B32").Value = 33.67373275756836
$ws.Range("B33").Value = 30.0858268737793
$ws.Range("B34").Value = 36.17101287841797
$ws.Range("B35").Value = 33.23563385009766
$ws.Range("B36").Value = 29.02166748046875
